$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2014-10"
$ws.Range("B2").Value = 104.6551
$ws.Range("C2").Value = 99.07729999999999
$ws.Range("D2").Value = 101.9783
$ws.Range("E2").Value = 101.5856
$ws.Range("F2").Value = 102.6077
$ws.Range("G2").Value = 101.6557
$ws.Range("H2").Value = 100.5971

$ws.Range("A3").Value = "2014-11"
$ws.Range("B3").Value = 102.7599
$ws.Range("C3").Value = 98.6662
$ws.Range("D3").Value = 101.7525
$ws.Range("E3").Value = 101.6731
$ws.Range("F3").Value = 102.4161
$ws.Range("G3").Value = 101.1482
$ws.Range("H3").Value = 100.3358

$ws.Range("A4").Value = "2014-12"
$ws.Range("B4").Value = 101.0155
$ws.Range("C4").Value = 98.4482
$ws.Range("D4").Value = 101.554
$ws.Range("E4").Value = 101.7847
$ws.Range("F4").Value = 102.7073
$ws.Range("G4").Value = 101.2641
$ws.Range("H4").Value = 100.3895

$ws.Range("A5").Value = "2014-01"
$ws.Range("B5").Value = 111.4854
$ws.Range("C5").Value = 99.0035
$ws.Range("D5").Value = 102.0521
$ws.Range("E5").Value = 100.7328
$ws.Range("F5").Value = 99.80929999999999
$ws.Range("G5").Value = 101.4322
$ws.Range("H5").Value = 101.5899

$ws.Range("A6").Value = "2014-02"
$ws.Range("B6").Value = 112.0981
$ws.Range("C6").Value = 99.0767
$ws.Range("D6").Value = 102.1306
$ws.Range("E6").Value = 100.8951
$ws.Range("F6").Value = 100.1635
$ws.Range("G6").Value = 101.5291
$ws.Range("H6").Value = 101.5038

$ws.Range("A7").Value = "2014-03"
$ws.Range("B7").Value = 111.8694
$ws.Range("C7").Value = 98.58620000000001
$ws.Range("D7").Value = 102.0722
$ws.Range("E7").Value = 100.9569
$ws.Range("F7").Value = 100.3375
$ws.Range("G7").Value = 101.3846
$ws.Range("H7").Value = 101.3138

$ws.Range("A8").Value = "2014-04"
$ws.Range("B8").Value = 111.2474
$ws.Range("C8").Value = 99.0488
$ws.Range("D8").Value = 102.2142
$ws.Range("E8").Value = 101.0132
$ws.Range("F8").Value = 100.2396
$ws.Range("G8").Value = 101.3945
$ws.Range("H8").Value = 101.1977

$ws.Range("A9").Value = "2014-05"
$ws.Range("B9").Value = 110.3402
$ws.Range("C9").Value = 99.2253
$ws.Range("D9").Value = 102.3928
$ws.Range("E9").Value = 101.0101
$ws.Range("F9").Value = 100.3108
$ws.Range("G9").Value = 101.3295
$ws.Range("H9").Value = 101.2286

$ws.Range("A10").Value = "2014-06"
$ws.Range("B10").Value = 109.5875
$ws.Range("C10").Value = 99.4808
$ws.Range("D10").Value = 102.2887
$ws.Range("E10").Value = 101.1714
$ws.Range("F10").Value = 100.7266
$ws.Range("G10").Value = 101.6899
$ws.Range("H10").Value = 101.2765

$ws.Range("A11").Value = "2014-07"
$ws.Range("B11").Value = 108.9108
$ws.Range("C11").Value = 98.8271
$ws.Range("D11").Value = 102.3205
$ws.Range("E11").Value = 101.0888
$ws.Range("F11").Value = 101.4125
$ws.Range("G11").Value = 101.6951
$ws.Range("H11").Value = 101.3032

$ws.Range("A12").Value = "2014-08"
$ws.Range("B12").Value = 108.0499
$ws.Range("C12").Value = 99.28660000000001
$ws.Range("D12").Value = 101.8861
$ws.Range("E12").Value = 101.2342
$ws.Range("F12").Value = 101.6507
$ws.Range("G12").Value = 101.6509
$ws.Range("H12").Value = 101.0838

$ws.Range("A13").Value = "2014-09"
$ws.Range("B13").Value = 105.7213
$ws.Range("C13").Value = 98.98699999999999
$ws.Range("D13").Value = 102.0033
$ws.Range("E13").Value = 101.4628
$ws.Range("F13").Value = 101.9036
$ws.Range("G13").Value = 101.3376
$ws.Range("H13").Value = 100.7607

$ws.Range("A14").Value = "2015-10"
$ws.Range("B14").Value = 97.5
$ws.Range("C14").Value = 98.90000000000001
$ws.Range("D14").Value = 100.9
$ws.Range("E14").Value = 100.7
$ws.Range("F14").Value = 101.3
$ws.Range("G14").Value = 99.40000000000001
$ws.Range("H14").Value = 100.6

$ws.Range("A15").Value = "2015-11"
$ws.Range("B15").Value = 97.8693
$ws.Range("C15").Value = 98.6404
$ws.Range("D15").Value = 100.8601
$ws.Range("E15").Value = 100.5114
$ws.Range("F15").Value = 101.1007
$ws.Range("G15").Value = 99.2303
$ws.Range("H15").Value = 100.6459

$ws.Range("A16").Value = "2015-12"
$ws.Range("B16").Value = 98.37779999999999
$ws.Range("C16").Value = 99.0068
$ws.Range("D16").Value = 100.9703
$ws.Range("E16").Value = 100.2688
$ws.Range("F16").Value = 100.8099
$ws.Range("G16").Value = 98.8034
$ws.Range("H16").Value = 100.4598

$ws.Range("A17").Value = "2015-01"
$ws.Range("B17").Value = 99.2801
$ws.Range("C17").Value = 98.8325
$ws.Range("D17").Value = 101.5809
$ws.Range("E17").Value = 101.5156
$ws.Range("F17").Value = 103.0294
$ws.Range("G17").Value = 100.8726
$ws.Range("H17").Value = 100.7294

$ws.Range("A18").Value = "2015-02"
$ws.Range("B18").Value = 97.61239999999999
$ws.Range("C18").Value = 98.7615
$ws.Range("D18").Value = 101.2512
$ws.Range("E18").Value = 101.3572
$ws.Range("F18").Value = 103.3585
$ws.Range("G18").Value = 100.9564
$ws.Range("H18").Value = 100.6331

$ws.Range("A19").Value = "2015-03"
$ws.Range("B19").Value = 97.283
$ws.Range("C19").Value = 98.8938
$ws.Range("D19").Value = 101.1333
$ws.Range("E19").Value = 101.2602
$ws.Range("F19").Value = 102.9068
$ws.Range("G19").Value = 100.9366
$ws.Range("H19").Value = 100.4371

$ws.Range("A20").Value = "2015-04"
$ws.Range("B20").Value = 96.9954
$ws.Range("C20").Value = 98.7582
$ws.Range("D20").Value = 101.2289
$ws.Range("E20").Value = 101.0958
$ws.Range("F20").Value = 102.9847
$ws.Range("G20").Value = 100.8146
$ws.Range("H20").Value = 100.659

$ws.Range("A21").Value = "2015-05"
$ws.Range("B21").Value = 97.00239999999999
$ws.Range("C21").Value = 99.40430000000001
$ws.Range("D21").Value = 101.3273
$ws.Range("E21").Value = 101.182
$ws.Range("F21").Value = 103.0455
$ws.Range("G21").Value = 100.6725
$ws.Range("H21").Value = 100.7828

$ws.Range("A22").Value = "2015-06"
$ws.Range("B22").Value = 97.21810000000001
$ws.Range("C22").Value = 99.2338
$ws.Range("D22").Value = 101.2733
$ws.Range("E22").Value = 101.0047
$ws.Range("F22").Value = 102.7105
$ws.Range("G22").Value = 100.4485
$ws.Range("H22").Value = 100.7399

$ws.Range("A23").Value = "2015-07"
$ws.Range("B23").Value = 97.2495
$ws.Range("C23").Value = 99.08159999999999
$ws.Range("D23").Value = 101.3871
$ws.Range("E23").Value = 101.1789
$ws.Range("F23").Value = 102.4973
$ws.Range("G23").Value = 100.1378
$ws.Range("H23").Value = 100.6227

$ws.Range("A24").Value = "2015-08"
$ws.Range("B24").Value = 97.496
$ws.Range("C24").Value = 99.54170000000001
$ws.Range("D24").Value = 101.4453
$ws.Range("E24").Value = 100.9866
$ws.Range("F24").Value = 102.2697
$ws.Range("G24").Value = 99.84010000000001
$ws.Range("H24").Value = 100.5901

$ws.Range("A25").Value = "2015-09"
$ws.Range("B25").Value = 97.4622
$ws.Range("C25").Value = 98.952
$ws.Range("D25").Value = 101.2986
$ws.Range("E25").Value = 100.651
$ws.Range("F25").Value = 101.8722
$ws.Range("G25").Value = 99.76390000000001
$ws.Range("H25").Value = 100.7347

$ws.Range("A26").Value = "2016-10"
$ws.Range("B26").Value = 98.90000000000001
$ws.Range("C26").Value = 101.3
$ws.Range("D26").Value = 100.2
$ws.Range("E26").Value = 99.8
$ws.Range("F26").Value = 100
$ws.Range("G26").Value = 99.7
$ws.Range("H26").Value = 99.40000000000001

$ws.Range("A27").Value = "2016-11"
$ws.Range("B27").Value = 98.8
$ws.Range("C27").Value = 101.7
$ws.Range("D27").Value = 100.4
$ws.Range("E27").Value = 99.7
$ws.Range("F27").Value = 99.90000000000001
$ws.Range("G27").Value = 99.59999999999999
$ws.Range("H27").Value = 99.5

$ws.Range("A28").Value = "2016-12"
$ws.Range("B28").Value = 99
$ws.Range("C28").Value = 102.7
$ws.Range("D28").Value = 100.8
$ws.Range("E28").Value = 99.7
$ws.Range("F28").Value = 100.1
$ws.Range("G28").Value = 99.8
$ws.Range("H28").Value = 99.8

$ws.Range("A29").Value = "2016-01"
$ws.Range("B29").Value = 98.851
$ws.Range("C29").Value = 99.0504
$ws.Range("D29").Value = 100.8428
$ws.Range("E29").Value = 99.9161
$ws.Range("F29").Value = 100.3676
$ws.Range("G29").Value = 99.67829999999999
$ws.Range("H29").Value = 100.3806

$ws.Range("A30").Value = "2016-02"
$ws.Range("B30").Value = 99.0487
$ws.Range("C30").Value = 98.8935
$ws.Range("D30").Value = 100.7885
$ws.Range("E30").Value = 99.967
$ws.Range("F30").Value = 99.4147
$ws.Range("G30").Value = 99.2377
$ws.Range("H30").Value = 100.3603

$ws.Range("A31").Value = "2016-03"
$ws.Range("B31").Value = 98.9623
$ws.Range("C31").Value = 99.2758
$ws.Range("D31").Value = 100.7739
$ws.Range("E31").Value = 100.0206
$ws.Range("F31").Value = 99.622
$ws.Range("G31").Value = 99.5942
$ws.Range("H31").Value = 100.302

$ws.Range("A32").Value = "2016-04"
$ws.Range("B32").Value = 98.8353
$ws.Range("C32").Value = 99.1519
$ws.Range("D32").Value = 100.7912
$ws.Range("E32").Value = 99.9995
$ws.Range("F32").Value = 99.9584
$ws.Range("G32").Value = 99.0813
$ws.Range("H32").Value = 99.7764

$ws.Range("A33").Value = "2016-05"
$ws.Range("B33").Value = 98.8
$ws.Range("C33").Value = 99.09999999999999
$ws.Range("D33").Value = 100.8
$ws.Range("E33").Value = 99.90000000000001
$ws.Range("F33").Value = 99.3
$ws.Range("G33").Value = 99
$ws.Range("H33").Value = 99.40000000000001

$ws.Range("A34").Value = "2016-06"
$ws.Range("B34").Value = 98.3
$ws.Range("C34").Value = 99.7
$ws.Range("D34").Value = 100.7
$ws.Range("E34").Value = 100.2
$ws.Range("F34").Value = 99.59999999999999
$ws.Range("G34").Value = 99.09999999999999
$ws.Range("H34").Value = 99.59999999999999

$ws.Range("A35").Value = "2016-07"
$ws.Range("B35").Value = 97.8
$ws.Range("C35").Value = 100.8
$ws.Range("D35").Value = 100.4
$ws.Range("E35").Value = 99.40000000000001
$ws.Range("F35").Value = 99.5
$ws.Range("G35").Value = 99.09999999999999
$ws.Range("H35").Value = 99.40000000000001

$ws.Range("A36").Value = "2016-08"
$ws.Range("B36").Value = 97.90000000000001
$ws.Range("C36").Value = 100.7
$ws.Range("D36").Value = 100.2
$ws.Range("E36").Value = 99.7
$ws.Range("F36").Value = 99.59999999999999
$ws.Range("G36").Value = 99.5
$ws.Range("H36").Value = 99.3

$ws.Range("A37").Value = "2016-09"
$ws.Range("B37").Value = 98.7
$ws.Range("C37").Value = 100.9
$ws.Range("D37").Value = 100.1
$ws.Range("E37").Value = 99.90000000000001
$ws.Range("F37").Value = 99.7
$ws.Range("G37").Value = 99.59999999999999
$ws.Range("H37").Value = 99.3

$ws.Range("A38").Value = "2017-10"
$ws.Range("B38").Value = 101.5
$ws.Range("C38").Value = 102
$ws.Range("D38").Value = 100.9
$ws.Range("E38").Value = 101.7
$ws.Range("F38").Value = 100.7
$ws.Range("G38").Value = 100.1
$ws.Range("H38").Value = 101.5

$ws.Range("A39").Value = "2017-11"
$ws.Range("B39").Value = 101.4
$ws.Range("C39").Value = 101.9
$ws.Range("D39").Value = 100.9
$ws.Range("E39").Value = 101.8
$ws.Range("F39").Value = 101
$ws.Range("G39").Value = 100.2
$ws.Range("H39").Value = 101.6

$ws.Range("A40").Value = "2017-12"
$ws.Range("B40").Value = 101.4
$ws.Range("C40").Value = 101.4
$ws.Range("D40").Value = 100.5
$ws.Range("E40").Value = 102.2
$ws.Range("F40").Value = 100.9
$ws.Range("G40").Value = 100.2
$ws.Range("H40").Value = 101.8

$ws.Range("A41").Value = "2017-01"
$ws.Range("B41").Value = 99.40000000000001
$ws.Range("C41").Value = 103
$ws.Range("D41").Value = 100.4
$ws.Range("E41").Value = 100
$ws.Range("F41").Value = 100.6
$ws.Range("G41").Value = 100
$ws.Range("H41").Value = 100.2

$ws.Range("A42").Value = "2017-02"
$ws.Range("B42").Value = 99.8
$ws.Range("C42").Value = 102.9
$ws.Range("D42").Value = 100.6
$ws.Range("E42").Value = 100.2
$ws.Range("F42").Value = 101.5
$ws.Range("G42").Value = 100.2
$ws.Range("H42").Value = 100.6

$ws.Range("A43").Value = "2017-03"
$ws.Range("B43").Value = 100.2
$ws.Range("C43").Value = 102.6
$ws.Range("D43").Value = 100.5
$ws.Range("E43").Value = 100.3
$ws.Range("F43").Value = 101.5
$ws.Range("G43").Value = 100.1
$ws.Range("H43").Value = 100.8

$ws.Range("A44").Value = "2017-04"
$ws.Range("B44").Value = 100.3
$ws.Range("C44").Value = 103
$ws.Range("D44").Value = 100.5
$ws.Range("E44").Value = 100.5
$ws.Range("F44").Value = 101.5
$ws.Range("G44").Value = 100.1
$ws.Range("H44").Value = 101.5

$ws.Range("A45").Value = "2017-05"
$ws.Range("B45").Value = 100.7
$ws.Range("C45").Value = 103.1
$ws.Range("D45").Value = 100.3
$ws.Range("E45").Value = 100.6
$ws.Range("F45").Value = 102
$ws.Range("G45").Value = 100.5
$ws.Range("H45").Value = 101.5

$ws.Range("A46").Value = "2017-06"
$ws.Range("B46").Value = 101.2
$ws.Range("C46").Value = 102.4
$ws.Range("D46").Value = 100.5
$ws.Range("E46").Value = 100.7
$ws.Range("F46").Value = 101.7
$ws.Range("G46").Value = 100.4
$ws.Range("H46").Value = 101.2

$ws.Range("A47").Value = "2017-07"
$ws.Range("B47").Value = 101.7
$ws.Range("C47").Value = 101.5
$ws.Range("D47").Value = 100.8
$ws.Range("E47").Value = 101.3
$ws.Range("F47").Value = 102.1
$ws.Range("G47").Value = 100.7
$ws.Range("H47").Value = 101

$ws.Range("A48").Value = "2017-08"
$ws.Range("B48").Value = 101.8
$ws.Range("C48").Value = 101.8
$ws.Range("D48").Value = 100.9
$ws.Range("E48").Value = 101.6
$ws.Range("F48").Value = 102
$ws.Range("G48").Value = 100.4
$ws.Range("H48").Value = 100.9

$ws.Range("A49").Value = "2017-09"
$ws.Range("B49").Value = 101.5
$ws.Range("C49").Value = 101.8
$ws.Range("D49").Value = 101
$ws.Range("E49").Value = 101.7
$ws.Range("F49").Value = 101.4
$ws.Range("G49").Value = 100.3
$ws.Range("H49").Value = 101.1
